$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove cell H2 entirely (column H data dropped for this row in the new export)
$ws.Range("H2").ClearContents()

# Updated simulation results for rows 2-25 (columns B:G, I:N)
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.097205796947011
$ws.Range("D2").Value = 1.099856455774556
$ws.Range("E2").Value = 1.095262943555872
$ws.Range("F2").Value = 1.100346054707249
$ws.Range("G2").Value = 1
$ws.Range("I2").Value = 1.03003759695652
$ws.Range("J2").Value = 1.102003611666362
$ws.Range("K2").Value = 1.102480098163674
$ws.Range("L2").Value = 1.097898192661956
$ws.Range("M2").Value = 1.102968466152379
$ws.Range("N2").Value = 1.103568583196509
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.102241083106902
$ws.Range("D3").Value = 1.104814097195892
$ws.Range("E3").Value = 1.099729372610027
$ws.Range("F3").Value = 1.10477024681982
$ws.Range("G3").Value = 1
$ws.Range("I3").Value = 1.03024986411576
$ws.Range("J3").Value = 1.106688352599673
$ws.Range("K3").Value = 1.107249946084999
$ws.Range("L3").Value = 1.102177057440727
$ws.Range("M3").Value = 1.107206197282106
$ws.Range("N3").Value = 1.108259976999293
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.105455937361601
$ws.Range("D4").Value = 1.107979333719769
$ws.Range("E4").Value = 1.102578944593439
$ws.Range("F4").Value = 1.107592016347746
$ws.Range("G4").Value = 1
$ws.Range("I4").Value = 1.030381538642308
$ws.Range("J4").Value = 1.109677473699487
$ws.Range("K4").Value = 1.110293935159551
$ws.Range("L4").Value = 1.104905431633896
$ws.Range("M4").Value = 1.109907466102745
$ws.Range("N4").Value = 1.111253342993926
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.106797496022318
$ws.Range("D5").Value = 1.109300174621315
$ws.Range("E5").Value = 1.103767568564554
$ws.Range("F5").Value = 1.108768834415345
$ws.Range("G5").Value = 1
$ws.Range("I5").Value = 1.030435559644966
$ws.Range("J5").Value = 1.110924373127348
$ws.Range("K5").Value = 1.111563857642513
$ws.Range("L5").Value = 1.106043135490101
$ws.Range("M5").Value = 1.111033653248121
$ws.Range("N5").Value = 1.112502013161995
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.107022176622148
$ws.Range("D6").Value = 1.109521384698505
$ws.Range("E6").Value = 1.103966606514093
$ws.Range("F6").Value = 1.10896588306443
$ws.Range("G6").Value = 1
$ws.Range("I6").Value = 1.030444552404358
$ws.Range("J6").Value = 1.111133173396939
$ws.Range("K6").Value = 1.111776521175889
$ws.Range("L6").Value = 1.106233625316822
$ws.Range("M6").Value = 1.111222202109348
$ws.Range("N6").Value = 1.112711109951919
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.10547390197165
$ws.Range("D7").Value = 1.107997020951994
$ws.Range("E7").Value = 1.102594863261474
$ws.Range("F7").Value = 1.107607777733405
$ws.Range("G7").Value = 1
$ws.Range("I7").Value = 1.030382265687932
$ws.Range("J7").Value = 1.10969417255856
$ws.Range("K7").Value = 1.110310941818156
$ws.Range("L7").Value = 1.104920669792995
$ws.Range("M7").Value = 1.109922550858633
$ws.Range("N7").Value = 1.111270065567294
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.098916711675844
$ws.Range("D8").Value = 1.101540996199558
$ws.Range("E8").Value = 1.096781003533197
$ws.Range("F8").Value = 1.101849935434423
$ws.Range("G8").Value = 1
$ws.Range("I8").Value = 1.030110523787943
$ws.Range("J8").Value = 1.103595817710699
$ws.Range("K8").Value = 1.104101111340473
$ws.Range("L8").Value = 1.099352823934978
$ws.Range("M8").Value = 1.104409294624195
$ws.Range("N8").Value = 1.105163050356058
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.087011415254083
$ws.Range("D9").Value = 1.089819108464829
$ws.Range("E9").Value = 1.086209014071905
$ws.Range("F9").Value = 1.091373258420965
$ws.Range("G9").Value = 1
$ws.Range("I9").Value = 1.029587071389653
$ws.Range("J9").Value = 1.092508486732914
$ws.Range("K9").Value = 1.09281556913983
$ws.Range("L9").Value = 1.089216056672642
$ws.Range("M9").Value = 1.094365187006086
$ws.Range("N9").Value = 1.09405997409655
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.078810870145601
$ws.Range("D10").Value = 1.081744832987181
$ws.Range("E10").Value = 1.078916015947601
$ws.Range("F10").Value = 1.084141841622692
$ws.Range("G10").Value = 1
$ws.Range("I10").Value = 1.029206438559473
$ws.Range("J10").Value = 1.084861167543222
$ws.Range("K10").Value = 1.08503457619634
$ws.Range("L10").Value = 1.082214976589448
$ws.Range("M10").Value = 1.087423815644835
$ws.Range("N10").Value = 1.086401794836449
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.075190749844811
$ws.Range("D11").Value = 1.078180466214283
$ws.Range("E11").Value = 1.07569399550006
$ws.Range("F11").Value = 1.080946099028364
$ws.Range("G11").Value = 1
$ws.Range("I11").Value = 1.029033709559118
$ws.Range("J11").Value = 1.081482851981503
$ws.Range("K11").Value = 1.08159793430711
$ws.Range("L11").Value = 1.079119926228479
$ws.Range("M11").Value = 1.084354204617634
$ws.Range("N11").Value = 1.083018681679134
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.073835070191877
$ws.Range("D12").Value = 1.076845669648886
$ws.Range("E12").Value = 1.074487020983404
$ws.Range("F12").Value = 1.079748834136754
$ws.Range("G12").Value = 1
$ws.Range("I12").Value = 1.028968325018433
$ws.Range("J12").Value = 1.080217362083851
$ws.Range("K12").Value = 1.080310705292032
$ws.Range("L12").Value = 1.077960212611094
$ws.Range("M12").Value = 1.083203884726891
$ws.Range("N12").Value = 1.081751394640673
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.074126376227766
$ws.Range("D13").Value = 1.077132488131576
$ws.Range("E13").Value = 1.074746390507556
$ws.Range("F13").Value = 1.080006123109426
$ws.Range("G13").Value = 1
$ws.Range("I13").Value = 1.028982406283243
$ws.Range("J13").Value = 1.080489304675743
$ws.Range("K13").Value = 1.080587314474829
$ws.Range("L13").Value = 1.078209439830369
$ws.Range("M13").Value = 1.083451099361367
$ws.Range("N13").Value = 1.082023723422238
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.075078917081908
$ws.Range("D14").Value = 1.078070356038701
$ws.Range("E14").Value = 1.075594437376597
$ws.Range("F14").Value = 1.080847344506427
$ws.Range("G14").Value = 1
$ws.Range("I14").Value = 1.029028330067507
$ws.Range("J14").Value = 1.081378466524621
$ws.Range("K14").Value = 1.081491753443892
$ws.Range("L14").Value = 1.079024272615704
$ws.Range("M14").Value = 1.084259328586282
$ws.Range("N14").Value = 1.082914147982931
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.075664331863487
$ws.Range("D15").Value = 1.078646753721754
$ws.Range("E15").Value = 1.076115582312352
$ws.Range("F15").Value = 1.081364277485475
$ws.Range("G15").Value = 1
$ws.Range("I15").Value = 1.029056461731918
$ws.Range("J15").Value = 1.081924881849882
$ws.Range("K15").Value = 1.082047571459123
$ws.Range("L15").Value = 1.07952496672795
$ws.Range("M15").Value = 1.084755946872524
$ws.Range("N15").Value = 1.083461339280627
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.079049616210123
$ws.Range("D16").Value = 1.081979902574649
$ws.Range("E16").Value = 1.079128454725153
$ws.Range("F16").Value = 1.084352528857266
$ws.Range("G16").Value = 1
$ws.Range("I16").Value = 1.029217732239269
$ws.Range("J16").Value = 1.085083916176196
$ws.Range("K16").Value = 1.085261185771996
$ws.Range("L16").Value = 1.082419001855718
$ws.Range("M16").Value = 1.087626144411725
$ws.Range("N16").Value = 1.086624859798031
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.08115416507905
$ws.Range("D17").Value = 1.084052044917847
$ws.Range("E17").Value = 1.081000815085156
$ws.Range("F17").Value = 1.086209346237634
$ws.Range("G17").Value = 1
$ws.Range("I17").Value = 1.029316749470137
$ws.Range("J17").Value = 1.087047171298184
$ws.Range("K17").Value = 1.087258553390899
$ws.Range("L17").Value = 1.084216979833647
$ws.Range("M17").Value = 1.089409063251896
$ws.Range("N17").Value = 1.088590902967481
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.082375075823629
$ws.Range("D18").Value = 1.085254156655257
$ws.Range("E18").Value = 1.082086784334738
$ws.Range("F18").Value = 1.087286211665748
$ws.Range("G18").Value = 1
$ws.Range("I18").Value = 1.02937374368968
$ws.Range("J18").Value = 1.088185882890261
$ws.Range("K18").Value = 1.088417119833208
$ws.Range("L18").Value = 1.085259616932336
$ws.Range("M18").Value = 1.090442876096128
$ws.Range("N18").Value = 1.08973123166063
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.082790266859505
$ws.Range("D19").Value = 1.085662954869026
$ws.Range("E19").Value = 1.082456045092965
$ws.Range("F19").Value = 1.087652361627115
$ws.Range("G19").Value = 1
$ws.Range("I19").Value = 1.029393049264512
$ws.Range("J19").Value = 1.088573081624227
$ws.Range("K19").Value = 1.088811081582456
$ws.Range("L19").Value = 1.085614111287583
$ws.Range("M19").Value = 1.090794354538268
$ws.Range("N19").Value = 1.090118980261211
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.080929057654547
$ws.Range("D20").Value = 1.083830403652535
$ws.Range("E20").Value = 1.080800568150406
$ws.Range("F20").Value = 1.086010770784566
$ws.Range("G20").Value = 1
$ws.Range("I20").Value = 1.029306204814299
$ws.Range("J20").Value = 1.08683720086743
$ws.Range("K20").Value = 1.087044927448315
$ws.Range("L20").Value = 1.08402470776275
$ws.Range("M20").Value = 1.089218411027042
$ws.Range("N20").Value = 1.088380634354631
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.07479872639588
$ws.Range("D21").Value = 1.077794481275123
$ws.Range("E21").Value = 1.075344994040854
$ws.Range("F21").Value = 1.08059991245492
$ws.Range("G21").Value = 1
$ws.Range("I21").Value = 1.029014840786671
$ws.Range("J21").Value = 1.081116928803563
$ws.Range("K21").Value = 1.081225719098771
$ws.Range("L21").Value = 1.078784607168648
$ws.Range("M21").Value = 1.084021609200794
$ws.Range("N21").Value = 1.082652238848309
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.070880341159736
$ws.Range("D22").Value = 1.073936466285643
$ws.Range("E22").Value = 1.071855709703944
$ws.Range("F22").Value = 1.077138450654387
$ws.Range("G22").Value = 1
$ws.Range("I22").Value = 1.028824539695607
$ws.Range("J22").Value = 1.077458538841169
$ws.Range("K22").Value = 1.077504691170835
$ws.Range("L22").Value = 1.075431376945547
$ws.Range("M22").Value = 1.080695284768548
$ws.Range("N22").Value = 1.078988653552545
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.072963829640224
$ws.Range("D23").Value = 1.07598785132572
$ws.Range("E23").Value = 1.073711241477137
$ws.Range("F23").Value = 1.078979258383869
$ws.Range("G23").Value = 1
$ws.Range("I23").Value = 1.02892610854825
$ws.Range("J23").Value = 1.079403980845606
$ws.Range("K23").Value = 1.079483382327575
$ws.Range("L23").Value = 1.07721472449353
$ws.Range("M23").Value = 1.082464396802666
$ws.Range("N23").Value = 1.080936858307774
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.081030794477033
$ws.Range("D24").Value = 1.083930573948909
$ws.Range("E24").Value = 1.080891070077509
$ws.Range("F24").Value = 1.086100517558573
$ws.Range("G24").Value = 1
$ws.Range("I24").Value = 1.029310971838989
$ws.Range("J24").Value = 1.086932097271684
$ws.Range("K24").Value = 1.087141475751509
$ws.Range("L24").Value = 1.084111606026697
$ws.Range("M24").Value = 1.089304577477982
$ws.Range("N24").Value = 1.088475665522664
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.090133556350051
$ws.Range("D25").Value = 1.092893170982293
$ws.Range("E25").Value = 1.088983394145656
$ws.Range("F25").Value = 1.094123368340769
$ws.Range("G25").Value = 1
$ws.Range("I25").Value = 1.029727849031466
$ws.Range("J25").Value = 1.092508486732914
$ws.Range("K25").Value = 1.095776458364883
$ws.Range("L25").Value = 1.091877656641324
$ws.Range("M25").Value = 1.097003220035098
$ws.Range("N25").Value = 1.096973508376485

Write-Host "Updated vm_pu.xlsx values for Case_0_145 (380 kV case)"
